# PHOENIX-5999: Modified the user of assistant engineer
#
# The "engineer" / "engineer1" rows (and their mirrored "assis_Engineer" /
# "assis_Engineer_1" rows) had their short dataName keys pointing at the
# wrong approver. Swap the dataName (column A) values between the paired
# rows so each key now lines up with the other approver's row, while every
# other column (department/designation/approver/remarks) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

# Row 7 <-> Row 8 (engineer / engineer1)
$ws.Range("A7").Value = "engineer1"
$ws.Range("A8").Value = "engineer"

# Row 19 <-> Row 20 (assis_Engineer / assis_Engineer_1)
$ws.Range("A19").Value = "assis_Engineer_1"
$ws.Range("A20").Value = "assis_Engineer"

# Reflect the cursor/selection state left behind in the saved file.
$excel.ActiveWindow.TabRatio = 0.5
$null = $ws.Range("A8").Select()
